# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# report sheets, as part of re-generating the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 17:17:18"
$wsZhCn.Range("H2").Value = "2016-03-22 17:17:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 17:17:23"
$wsDeDe.Range("H2").Value = "2016-03-22 17:18:03"
